$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6691
$ws.Range("E2").Value = 634
$ws.Range("F2").Value = 634
$ws.Range("G2").Value = 644
$ws.Range("H2").Value = 457
$ws.Range("I2").Value = 498
$ws.Range("J2").Value = -41
$ws.Range("K2").Value = 8315
$ws.Range("L2").Value = 3019
$ws.Range("M2").Value = 5296
$ws.Range("N2").Value = 5254
$ws.Range("O2").Value = 41
$ws.Range("P2").Value = 329
$ws.Range("Q2").Value = 670
$ws.Range("R2").Value = -183
$ws.Range("S2").Value = -872
$ws.Range("T2").Value = 267
$ws.Range("U2").Value = 403
$ws.Range("V2").Value = 2121
$ws.Range("W2").Value = 9.48
$ws.Range("X2").Value = 6.83
$ws.Range("Y2").Value = 9.789999999999999
$ws.Range("Z2").Value = 5.34
$ws.Range("AA2").Value = 57.02
$ws.Range("AB2").Value = 1434.45
$ws.Range("AC2").Value = 5603
$ws.Range("AD2").Value = 7.61
$ws.Range("AE2").Value = 59104
$ws.Range("AF2").Value = 0.72
$ws.Range("AG2").Value = 815
$ws.Range("AH2").Value = 1.91
$ws.Range("AI2").Value = 14.54
$ws.Range("AJ2").Value = 8889843

# Row 3
$ws.Range("D3").Value = 7082
$ws.Range("E3").Value = 531
$ws.Range("F3").Value = 531
$ws.Range("G3").Value = 969
$ws.Range("H3").Value = 709
$ws.Range("I3").Value = 759
$ws.Range("J3").Value = -50
$ws.Range("K3").Value = 10368
$ws.Range("L3").Value = 4296
$ws.Range("M3").Value = 6072
$ws.Range("N3").Value = 5996
$ws.Range("O3").Value = 76
$ws.Range("P3").Value = 444
$ws.Range("Q3").Value = -36
$ws.Range("R3").Value = 255
$ws.Range("S3").Value = -369
$ws.Range("T3").Value = 246
$ws.Range("U3").Value = -282
$ws.Range("V3").Value = 3343
$ws.Range("W3").Value = 7.49
$ws.Range("X3").Value = 10.01
$ws.Range("Y3").Value = 13.49
$ws.Range("Z3").Value = 7.59
$ws.Range("AA3").Value = 70.76000000000001
$ws.Range("AB3").Value = 1187.99
$ws.Range("AC3").Value = 8536
$ws.Range("AD3").Value = 5.35
$ws.Range("AE3").Value = 67454
$ws.Range("AF3").Value = 0.68
$ws.Range("AG3").Value = 1100
$ws.Range("AH3").Value = 2.41
$ws.Range("AI3").Value = 12.89
$ws.Range("AJ3").Value = 8889843

# Row 4
$ws.Range("D4").Value = 7409
$ws.Range("E4").Value = 680
$ws.Range("F4").Value = 680
$ws.Range("G4").Value = 645
$ws.Range("H4").Value = 485
$ws.Range("I4").Value = 499
$ws.Range("J4").Value = -14
$ws.Range("K4").Value = 10468
$ws.Range("L4").Value = 4150
$ws.Range("M4").Value = 6317
$ws.Range("N4").Value = 6277
$ws.Range("O4").Value = 40
$ws.Range("P4").Value = 444
$ws.Range("Q4").Value = 1829
$ws.Range("R4").Value = -963
$ws.Range("S4").Value = -426
$ws.Range("T4").Value = 1403
$ws.Range("U4").Value = 426
$ws.Range("V4").Value = 3151
$ws.Range("W4").Value = 9.19
$ws.Range("X4").Value = 6.54
$ws.Range("Y4").Value = 8.130000000000001
$ws.Range("Z4").Value = 4.65
$ws.Range("AA4").Value = 65.7
$ws.Range("AB4").Value = 1272.16
$ws.Range("AC4").Value = 5609
$ws.Range("AD4").Value = 7.97
$ws.Range("AE4").Value = 70621
$ws.Range("AF4").Value = 0.63
$ws.Range("AG4").Value = 1100
$ws.Range("AH4").Value = 2.46
$ws.Range("AI4").Value = 19.61
$ws.Range("AJ4").Value = 8889843

# Row 5
$ws.Range("D5").Value = 7624
$ws.Range("E5").Value = 671
$ws.Range("F5").Value = 671
$ws.Range("G5").Value = 534
$ws.Range("H5").Value = 411
$ws.Range("I5").Value = 396
$ws.Range("J5").Value = 15
$ws.Range("K5").Value = 11131
$ws.Range("L5").Value = 4637
$ws.Range("M5").Value = 6495
$ws.Range("N5").Value = 6515
$ws.Range("O5").Value = -20
$ws.Range("P5").Value = 444
$ws.Range("Q5").Value = -190
$ws.Range("R5").Value = -899
$ws.Range("S5").Value = 776
$ws.Range("T5").Value = 1005
$ws.Range("U5").Value = -1195
$ws.Range("V5").Value = 3902
$ws.Range("W5").Value = 8.800000000000001
$ws.Range("X5").Value = 5.39
$ws.Range("Y5").Value = 6.19
$ws.Range("Z5").Value = 3.8
$ws.Range("AA5").Value = 71.39
$ws.Range("AB5").Value = 1338.51
$ws.Range("AC5").Value = 4454
$ws.Range("AD5").Value = 10.45
$ws.Range("AE5").Value = 73290
$ws.Range("AF5").Value = 0.64
$ws.Range("AG5").Value = 1100
$ws.Range("AH5").Value = 2.36
$ws.Range("AI5").Value = 24.7
$ws.Range("AJ5").Value = 8889843

# Row 6
$ws.Range("D6").Value = 8140
$ws.Range("E6").Value = 801
$ws.Range("F6").Value = 801
$ws.Range("G6").Value = 629
$ws.Range("H6").Value = 360
$ws.Range("I6").Value = 361
$ws.Range("K6").Value = 11119
$ws.Range("L6").Value = 4564
$ws.Range("M6").Value = 6554
$ws.Range("N6").Value = 6576
$ws.Range("P6").Value = 444
$ws.Range("Q6").Value = 550
$ws.Range("R6").Value = -483
$ws.Range("S6").Value = -238
$ws.Range("T6").Value = 732
$ws.Range("U6").Value = -183
$ws.Range("V6").Value = 3774
$ws.Range("W6").Value = 9.84
$ws.Range("X6").Value = 4.42
$ws.Range("Y6").Value = 5.52
$ws.Range("Z6").Value = 3.23
$ws.Range("AA6").Value = 69.64
$ws.Range("AB6").Value = 1468.23
$ws.Range("AC6").Value = 4062
$ws.Range("AD6").Value = 10.34
$ws.Range("AE6").Value = 73967
$ws.Range("AF6").Value = 0.57
$ws.Range("AG6").Value = 1100
$ws.Range("AH6").Value = 2.62
$ws.Range("AI6").Value = 27.08
$ws.Range("AJ6").Value = 8889843

# Row 7
$ws.Range("D7").Value = 9337
$ws.Range("E7").Value = 1047
$ws.Range("H7").Value = 658
$ws.Range("I7").Value = 660
$ws.Range("W7").Value = 11.21
$ws.Range("X7").Value = 7.05
$ws.Range("AC7").Value = 7424
$ws.Range("AD7").Value = 6.1
$ws.Range("G7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 10311
$ws.Range("E8").Value = 1096
$ws.Range("H8").Value = 791
$ws.Range("I8").Value = 793
$ws.Range("W8").Value = 10.63
$ws.Range("X8").Value = 7.67
$ws.Range("AC8").Value = 8920
$ws.Range("AD8").Value = 4.93
$ws.Range("G8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
